# Automatische test-sync: 2025-08-03 18:38:50
#
# Adds a new "Testmail #12" row to the Logs sheet and refreshes the
# Dashboard category-count table (counts + sort order) to account for it.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1. Logs sheet: append row 40
# ---------------------------------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$logs.Range("A40").Value = "Ik heb nog geen geld terug."
$logs.Range("B40").Value = "mailmind.test@zohomail.eu"
$logs.Range("C40").Value = "Testmail #12: Ik heb nog geen geld terug."
$logs.Range("D40").Value = "Retour / Terugbetaling"
$logs.Range("E40").Value = "Beste klant,`nDank u wel voor uw bericht. Om uw situatie te kunnen bekijken en u verder te kunnen helpen, hebben wij een aantal gegevens nodig. Kunt u ons alstublieft uw naam en het ordernummer doorgeven? Op die manier kunnen wij uw specifieke zaak onderzoeken en nagaan waarom het terugbetalingsproces nog niet is afgerond.`nWij kijken uit naar uw reactie.`nMet vriendelijke groet,`n[Bedrijfsnaam] E-mailassistent"
$logs.Range("F40").Value = "2025-08-03 18:38:04"
$logs.Range("G40").Value = "Ja"
$logs.Range("H40").Value = "Nee"
$logs.Range("I40").Value = "Ja"
$logs.Range("J40").Value = "Nee"

# ---------------------------------------------------------------
# 1b. Logs sheet: extend the conditional-formatting ranges so the
#     rules that used to stop at row 39 now also cover row 40.
# ---------------------------------------------------------------
foreach ($col in @("D", "G", "H", "I", "J")) {
    $oldRange = $logs.Range("$col" + "2:" + "$col" + "39")
    $newRange = $logs.Range("$col" + "2:" + "$col" + "40")
    for ($i = 1; $i -le $oldRange.FormatConditions.Count; $i++) {
        $oldRange.FormatConditions.Item($i).ModifyAppliesToRange($newRange)
    }
}

# ---------------------------------------------------------------
# 2. Dashboard sheet: update the "Aantal e-mails per categorie" table.
#    The new row bumps "Retour / Terugbetaling" from 1 -> 2, which
#    moves it to the top of the (count desc, then alphabetical) block
#    of categories, pushing the other count=1 categories down one row.
# ---------------------------------------------------------------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Range("A6").Value = "Retour / Terugbetaling"
$dash.Range("B6").Value = 2

$dash.Range("A7").Value = "Documentatie / Datasheets"
$dash.Range("B7").Value = 1

$dash.Range("A8").Value = "Klacht / Probleem"
$dash.Range("B8").Value = 1

$dash.Range("A9").Value = "Klantenservice / Contact"
$dash.Range("B9").Value = 1
